# This edit reorders the weekly price records in rows 2-14 (the header is
# row 1). The underlying data for each row (Mercado / Producto / Variedad /
# etc.) never changes - only which "Fecha" + price-quote row ends up on
# which line moves, consistent with a re-sort of the weekly rows.
#
# Columns that are touched per row: D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg), T (Kg/unidad)
# Columns A,B,C,E,F,G,H,I,J,K are identical on every row and stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by number) whose values move together as a record.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the "before" values for every data row (2-14) and every tracked
# column, so the row-to-row moves below don't clobber source data before it
# has been read.
$snapshot = @{}
for ($row = 2; $row -le 14; $row++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Cells.Item($row, $col).Value2
    }
    $snapshot[$row] = $rowVals
}

# Destination row -> source row (source row's record moves onto destination).
$rowMap = @{
    2  = 4
    3  = 5
    4  = 12
    5  = 11
    6  = 13
    7  = 6
    8  = 14
    9  = 2
    10 = 3
    11 = 9
    12 = 10
    13 = 7
    14 = 8
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value = $srcVals[$col]
    }
}

Write-Output "Reordered rows 2-14 per updated weekly sort"
